# Actualizacion 01-06-2020 carga masiva
# Rebuild the inventory template header row with the new Spanish field
# names, extend the sheet from 7 to 13 columns, and apply the new header
# look (bold white text on an accent-6 fill with thin borders).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Header captions (row 1, columns A:M)
# ---------------------------------------------------------------------
$ws.Range("A1").Value = "NOMBRE_NODO"
$ws.Range("B1").Value = "ID_TIPO"
$ws.Range("C1").Value = "SERIE"
$ws.Range("D1").Value = "DIRECCION_IP"
$ws.Range("E1").Value = "ID_REGION"
$ws.Range("F1").Value = "IOS_IMAGE"
$ws.Range("G1").Value = "IOS_VERSION"
$ws.Range("H1").Value = "LATITUD"
$ws.Range("I1").Value = "LONGITUD"
$ws.Range("J1").Value = "FECHA_MANTENIMIENTO"
$ws.Range("K1").Value = "ID_CONTRATO"
$ws.Range("L1").Value = "ID_PROVEEDOR"
$ws.Range("M1").Value = "ID_UBICACION"

# ---------------------------------------------------------------------
# 2. Header styling: bold white text on an accent-6 themed fill.
#    Columns A:K get a full thin border; L:M (the last two columns) get
#    a thin left/right-only border, matching the template's look.
# ---------------------------------------------------------------------
$headerMain = $ws.Range("A1:K1")
$headerMain.Font.Bold = $true
$headerMain.Font.ThemeColor = 2
$headerMain.Interior.Pattern = 1
$headerMain.Interior.ThemeColor = 7
$headerMain.Interior.TintAndShade = -0.249977111117893
$headerMain.Borders.LineStyle = 1
$headerMain.Borders.Weight = 2

$headerTail = $ws.Range("L1:M1")
$headerTail.Font.Bold = $true
$headerTail.Font.ThemeColor = 2
$headerTail.Interior.Pattern = 1
$headerTail.Interior.ThemeColor = 7
$headerTail.Interior.TintAndShade = -0.249977111117893

# Each header cell in this pair needs its own left+right thin border
# (not just the outer edges of the combined range), so set them one
# cell at a time.
foreach ($addr in @("L1", "M1")) {
    $cell = $ws.Range($addr)
    $cell.Borders.Item(7).LineStyle = 1
    $cell.Borders.Item(7).Weight = 2
    $cell.Borders.Item(10).LineStyle = 1
    $cell.Borders.Item(10).Weight = 2
}

# ---------------------------------------------------------------------
# 3. Column J (and the J1/K1 headers) carry a date number format, for
#    the maintenance-date values that belong under those headers.
# ---------------------------------------------------------------------
$ws.Range("J1").NumberFormat = "mm-dd-yy"
$ws.Range("K1").NumberFormat = "mm-dd-yy"
$ws.Columns.Item(10).NumberFormat = "mm-dd-yy"

# ---------------------------------------------------------------------
# 4. Column widths (characters), matching the published template.
# ---------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 14.6666666666667
$ws.Columns.Item(2).ColumnWidth = 7.16666666666667
$ws.Columns.Item(3).ColumnWidth = 4.83333333333333
$ws.Columns.Item(4).ColumnWidth = 12.6666666666667
$ws.Columns.Item(5).ColumnWidth = 10
$ws.Columns.Item(6).ColumnWidth = 10.3333333333333
$ws.Columns.Item(7).ColumnWidth = 12
$ws.Columns.Item(8).ColumnWidth = 7.66666666666667
$ws.Columns.Item(9).ColumnWidth = 9.66666666666667
$ws.Columns.Item(10).ColumnWidth = 23
$ws.Columns.Item(11).ColumnWidth = 12.8333333333333
$ws.Columns.Item(12).ColumnWidth = 13.8333333333333
$ws.Columns.Item(13).ColumnWidth = 13.1666666666667

# ---------------------------------------------------------------------
# 5. Selection mirrors the saved view state in the updated workbook.
# ---------------------------------------------------------------------
$ws.Range("K5").Select() | Out-Null
